# Generate Report for Handoff
# Updates the localization status report: flips Status from
# "In Translation" to "Ready for handoff" and refreshes the
# handoff/generate timestamps on all three sheets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet: status columns (zh-cn/de-de) + generate date ---
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 04:34:37"

# --- zh-cn sheet: status + latest handoff datetime ---
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-19 04:34:32"

# --- de-de sheet: status + latest handoff datetime ---
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-19 04:34:37"

# --- Column widths: Status columns widen to fit "Ready for handoff" ---
# (target stored width ~17.216 chars; engine quantizes ColumnWidth onto a
# 1/6-character grid, so 16.3 is the input that lands closest to it)
$wsOverview.Range("E:E").ColumnWidth = 16.3
$wsOverview.Range("F:F").ColumnWidth = 16.3
$wsZhCn.Range("C:C").ColumnWidth = 16.3
$wsDeDe.Range("C:C").ColumnWidth = 16.3
